# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (total) sheet,
#    populated with the per-fund holding breakdown for 2022-Q1.
# 2. Update the "总计" summary sheet with a new top data row for 2022-Q1,
#    pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")
$totalSheetBeforeInsert = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet (positioned before "总计")
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheetBeforeInsert, $null)
$newSheet.Name = "2022-Q1"

# NOTE: inserting a sheet shifts the position-based index of every sheet
# after it, which silently invalidates worksheet handles obtained before
# the insert. Re-resolve "总计" by name now that it has moved.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row text
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows: A (index), B (code, text), C (name, text), D (scale, text),
# E (total stock position, text), F (position ratio, text), G (held value, text), H (rank, number)
$data = @(
    @("0", "320003", "诺安先锋混合",             "45.79", "69.96", "3.04", "1.3920", 7),
    @("1", "010874", "泰康品质生活混合A",         "13.17", "81.43", "3.55", "0.4675", 4),
    @("2", "005014", "泰康景泰回报混合A",         "11.64", "27.21", "1.44", "0.1676", 6),
    @("3", "010875", "泰康品质生活混合C",         "4.39",  "81.43", "3.55", "0.1558", 4),
    @("4", "540007", "汇丰晋信中小盘股票",         "0.61",  "93.28", "4.23", "0.0258", 1),
    @("5", "005015", "泰康景泰回报混合C",         "0.63",  "27.21", "1.44", "0.0091", 6),
    @("6", "006143", "恒生前海中证质量成长低波动指数A", "0.06",  "94.34", "2.35", "0.0014", 7),
    @("7", "006144", "恒生前海中证质量成长低波动指数C", "0.01",  "94.34", "2.35", "0.0002", 7)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = [double]$row[0]

    $newSheet.Range("B$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("B$r").ClearFormats()

    $newSheet.Range("C$r").NumberFormat = "@"
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("C$r").ClearFormats()

    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("D$r").ClearFormats()

    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("E$r").ClearFormats()

    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("F$r").ClearFormats()

    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("G$r").ClearFormats()

    $newSheet.Cells.Item($r, 8).Value = [double]$row[7]

    $r = $r + 1
}

# Re-apply the bold/centered/bordered look (style used across every other
# sheet's header row + index column) by copying formatting from the
# "2021-Q4" template sheet, which keeps the same visual style.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert the 2022-Q1 total as the new first
#    data row, shifting the previously-existing rows down by one.
#    (Reading `.Value` back from a Range is unreliable on this host, so
#    the shifted rows are written from known literal values rather than
#    copied from the existing cells.)
# ---------------------------------------------------------------------
$totalRows = @(
    @(0, "2022-Q1", 8,  2.22),
    @(1, "2021-Q4", 16, 4.66),
    @(2, "2021-Q3", 4,  0.87),
    @(3, "2021-Q2", 1,  0.05),
    @(4, "2021-Q1", 5,  1.01)
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = [double]$row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = [double]$row[2]
    $totalSheet.Cells.Item($r, 4).Value = [double]$row[3]
    $r = $r + 1
}

# Row 6 is brand new (the sheet used to stop at row 5), so its "A" cell
# doesn't yet carry the bold/centered/bordered index-column look the
# other rows already have. Copy it over from an existing styled cell.
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
